$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 50, shifting the existing rows (and all rows below them) down by 2.
$ws.Range("A50:A51").EntireRow.Insert()

# Fill in the new row 50 (Primera quality, newest report date).
$ws.Cells.Item(50, 1).Value = 7
$ws.Cells.Item(50, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(50, 3).Value = "Ñuble"
$ws.Cells.Item(50, 4).Value = 44953
$ws.Cells.Item(50, 5).Value = 16
$ws.Cells.Item(50, 6).Value = 100112040
$ws.Cells.Item(50, 7).Value = "Cilantro"
$ws.Cells.Item(50, 8).Value = "Sin especificar"
$ws.Cells.Item(50, 9).Value = "Primera"
$ws.Cells.Item(50, 10).Value = 300
$ws.Cells.Item(50, 11).Value = 800
$ws.Cells.Item(50, 12).Value = 900
$ws.Cells.Item(50, 13).Value = 850
$ws.Cells.Item(50, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(50, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(50, 16).Value = 850
$ws.Cells.Item(50, 17).Value = 1
$ws.Cells.Item(50, 18).Value = "Hortaliza"

# Fill in the new row 51 (Segunda quality, newest report date).
$ws.Cells.Item(51, 1).Value = 7
$ws.Cells.Item(51, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(51, 3).Value = "Ñuble"
$ws.Cells.Item(51, 4).Value = 44953
$ws.Cells.Item(51, 5).Value = 16
$ws.Cells.Item(51, 6).Value = 100112040
$ws.Cells.Item(51, 7).Value = "Cilantro"
$ws.Cells.Item(51, 8).Value = "Sin especificar"
$ws.Cells.Item(51, 9).Value = "Segunda"
$ws.Cells.Item(51, 10).Value = 200
$ws.Cells.Item(51, 11).Value = 600
$ws.Cells.Item(51, 12).Value = 600
$ws.Cells.Item(51, 13).Value = 600
$ws.Cells.Item(51, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(51, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(51, 16).Value = 600
$ws.Cells.Item(51, 17).Value = 1
$ws.Cells.Item(51, 18).Value = "Hortaliza"
